$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.031.09'
$ws.Range("E2").Value = '  -0.45%  '

$ws.Range("D3").Value = '3.485.03'
$ws.Range("E3").Value = '  +4.92%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '648.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.47'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.415'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.82%  '

$ws.Range("D11").Value = '3.486.54'
$ws.Range("E11").Value = '  +5.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.66%  '

$ws.Range("E13").Value = '  -2.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").Value = '95.803.44'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '4.140.55'
$ws.Range("E16").Value = '  +5.00%  '

$ws.Range("E17").Value = '  +2.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("D19").Value = '3.480.19'
$ws.Range("E19").Value = '  +4.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +14.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.522'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '512.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000194'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '92.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.18%  '

$ws.Range("D29").Value = '3.665.11'
$ws.Range("E29").Value = '  +4.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.52%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.139'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.185'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '31.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.94%  '

$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.579'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.43%  '

$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.17%  '

$ws.Range("E39").Value = '  -2.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '522.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.77%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.151'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.35%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.923'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +12.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0421'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.72%  '

$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.15%  '

$ws.Range("B48").Value = 'MantraDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.50%  '

$ws.Range("E51").Value = '  +1.30%  '
